{"js": "// Apply the \"Added many more features\" edits to the Collapsed Castle\n// Bonus Buy review document.\n//\n// The title/meta string (\"Play Collapsed Castle Bonus Buy Free Slot -\n// Pros and Cons\") appears twice in the document (once as the H1 heading,\n// once as a bold run near the bottom) and both occurrences must be\n// updated to the same new text, so it is handled with a loop over all\n// search hits. Every other change is a single, unique text replacement.\n\nconst replacements = [\n  {\n    find: \"Play Collapsed Castle Bonus Buy Free Slot - Pros and Cons\",\n    replace: \"Play Collapsed Castle Bonus Buy | Free Slot Game\",\n  },\n  {\n    find: \"Original cascading reels\",\n    replace: \"Exciting gameplay with cascading reels\",\n  },\n  {\n    find: \"High payout potential\",\n    replace: \"Wide betting range suitable for all players\",\n  },\n  {\n    find: \"Beautifully designed medieval theme\",\n    replace: \"Impressive payout potential of up to 3,302x your bet\",\n  },\n  {\n    find: \"Suitable for both low rollers and high rollers\",\n    replace: \"Beautifully designed visuals and immersive sound\",\n  },\n  {\n    find: \"Requires payment to activate the Free Spins feature\",\n    replace: \"Limited number of free spins in the bonus game\",\n  },\n  {\n    find:\n      \"Read our review of Collapsed Castle Bonus Buy to discover the pros and cons of this medieval-themed slot machine. Play for free and see if it's the right game for you!\",\n    replace:\n      \"Read our review of Collapsed Castle Bonus Buy and play this free slot game with exciting gameplay and impressive payout potential.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Collapsed Castle\n# Bonus Buy review document via Word COM interop (Find/Replace).\n#\n# Note: the title/meta string (\"Play Collapsed Castle Bonus Buy Free\n# Slot - Pros and Cons\") appears twice (the H1 heading and a bold run\n# near the bottom) and both must become the same new text, so that one\n# uses wdReplaceAll (2). Every other string is unique in the document,\n# but MatchCase is kept on throughout because \"Original cascading\n# reels\" (the bullet, capitalized) must change while the unrelated,\n# lowercase \"original cascading reels\" phrase inside a body paragraph\n# must be left untouched.\n\n$d = $word.ActiveDocument\n$content = $d.Content\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"Play Collapsed Castle Bonus Buy Free Slot - Pros and Cons\" \"Play Collapsed Castle Bonus Buy | Free Slot Game\"\nReplace-Text \"Original cascading reels\" \"Exciting gameplay with cascading reels\"\nReplace-Text \"High payout potential\" \"Wide betting range suitable for all players\"\nReplace-Text \"Beautifully designed medieval theme\" \"Impressive payout potential of up to 3,302x your bet\"\nReplace-Text \"Suitable for both low rollers and high rollers\" \"Beautifully designed visuals and immersive sound\"\nReplace-Text \"Requires payment to activate the Free Spins feature\" \"Limited number of free spins in the bonus game\"\nReplace-Text \"Read our review of Collapsed Castle Bonus Buy to discover the pros and cons of this medieval-themed slot machine. Play for free and see if it's the right game for you!\" \"Read our review of Collapsed Castle Bonus Buy and play this free slot game with exciting gameplay and impressive payout potential.\"\n"}
